$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in 0 for previously-blank cells across several rows (mirrors the
# author filling out the "Сети" score sheet with 0 where no score existed).
$ws.Range("F4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = 0

$ws.Range("F5").Value = 0

$ws.Range("L6").Value = 0
$ws.Range("N6").Value = 0

$ws.Range("N7").Value = 0

$ws.Range("F8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("L8").Value = 0

$ws.Range("F9").Value = 0

$ws.Range("F14").Value = 0
$ws.Range("H14").Value = 0

$ws.Range("F16").Value = 0

$ws.Range("F19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = 0

$ws.Range("L20").Value = 0
$ws.Range("N20").Value = 0

$ws.Range("H21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = 0

$ws.Range("F24").Value = 0

# Update the active selection on the sheet to match the saved view state.
$ws.Range("N4:N25").Select()
